$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (A: bill_number, B: id_client, C: quantity, D: id_product)
$data = @(
    @(3, 2, 3, 3),
    @(4, 2, 2, 4),
    @(5, 2, 3, 6),
    @(6, 3, 1, 9),
    @(7, 4, 2, 5),
    @(8, 4, 3, 7),
    @(9, 4, 4, 2),
    @(10, 4, 1, 4)
)

$r = 4
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Formula = "=(D$r-1)*10000"
    $ws.Cells.Item($r, 6).Formula = "=C$r*E$r"
    $r = $r + 1
}

$ws.Range("G19").Select()
